$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Season Group" column (column B) entirely, shifting remaining
# columns (Y/P, Lng) left.
$ws.Range("B:B").Delete()

# Remove rows 3 and 4 (Group2 and Difference rows), leaving only the
# header row and the first data row.
$ws.Range("A3:A4").EntireRow.Delete()

# Update the remaining data row values per the new aggregate numbers.
$ws.Range("B2").Value = 4.4491525423729
$ws.Range("C2").Value = 5.555555555555562
